# Auto-generated Excel COM-interop edit script
# Applies the scheduled market-data refresh described in the commit diff
# (updates cached currentAveragePrice / Leve price / profit columns per sheet)

$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 34.4
$ws.Range("I6").Value = 34.4
$ws.Range("K6").Value = 103.2
$ws.Range("M6").Value = 8.800000000000011
$ws.Range("H9").Value = 2101102.8
$ws.Range("I9").Value = 3676741.5
$ws.Range("J9").Value = 251
$ws.Range("K9").Value = 3676741.5
$ws.Range("L9").Value = 251
$ws.Range("M9").Value = -3676572.5
$ws.Range("N9").Value = -589
$ws.Range("H12").Value = 2037.4286
$ws.Range("I12").Value = 422.66666
$ws.Range("J12").Value = 3248.5
$ws.Range("K12").Value = 422.66666
$ws.Range("L12").Value = 3248.5
$ws.Range("M12").Value = -252.66666
$ws.Range("N12").Value = -3588.5
$ws.Range("H32").Value = 16672426
$ws.Range("J32").Value = 8599.666999999999
$ws.Range("L32").Value = 8599.666999999999
$ws.Range("N32").Value = -9251.666999999999
$ws.Range("H38").Value = 50004.5
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H40").Value = 6350.185
$ws.Range("I40").Value = 5600.1113
$ws.Range("K40").Value = 5600.1113
$ws.Range("M40").Value = -5425.1113
$ws.Range("H43").Value = 999.5
$ws.Range("J43").Value = 999
$ws.Range("L43").Value = 999
$ws.Range("N43").Value = -1137
$ws.Range("H58").Value = 1179.5
$ws.Range("I58").Value = 1205.2858
$ws.Range("K58").Value = 3615.8574
$ws.Range("M58").Value = -3465.8574
$ws.Range("H94").Value = 3156.7144
$ws.Range("I94").Value = 1519.6
$ws.Range("K94").Value = 1519.6
$ws.Range("M94").Value = -1068.6
$ws.Range("H116").Value = 10455.111
$ws.Range("I116").Value = 5400
$ws.Range("J116").Value = 16774
$ws.Range("K116").Value = 5400
$ws.Range("L116").Value = 16774
$ws.Range("M116").Value = -1958
$ws.Range("N116").Value = -23658
$ws.Range("H131").Value = 2002097.4
$ws.Range("I131").Value = 2501497
$ws.Range("K131").Value = 7504491
$ws.Range("M131").Value = -7499451
$ws.Range("H135").Value = 1651.7778
$ws.Range("I135").Value = 1447.1428
$ws.Range("K135").Value = 13024.2852
$ws.Range("M135").Value = -10489.2852
$ws.Range("H137").Value = 1392453.1
$ws.Range("I137").Value = 1615782.4
$ws.Range("J137").Value = 7811.4
$ws.Range("K137").Value = 4847347.199999999
$ws.Range("L137").Value = 23434.2
$ws.Range("M137").Value = -4844797.199999999
$ws.Range("N137").Value = -28534.2
$ws.Range("H138").Value = 3124.4783
$ws.Range("J138").Value = 3222.1428
$ws.Range("L138").Value = 9666.428400000001
$ws.Range("N138").Value = -19946.4284
$ws.Range("H141").Value = 2733.4
$ws.Range("I141").Value = 2662.8572
$ws.Range("J141").Value = 2898
$ws.Range("K141").Value = 7988.571599999999
$ws.Range("L141").Value = 8694
$ws.Range("M141").Value = -2808.571599999999
$ws.Range("N141").Value = -19054

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2098.6155
$ws.Range("I2").Value = 1688.3
$ws.Range("J2").Value = 3466.3333
$ws.Range("K2").Value = 1688.3
$ws.Range("L2").Value = 3466.3333
$ws.Range("M2").Value = -1575.3
$ws.Range("N2").Value = -3692.3333
$ws.Range("H5").Value = 340.42856
$ws.Range("I5").Value = 314
$ws.Range("K5").Value = 314
$ws.Range("M5").Value = -202
$ws.Range("H32").Value = 1431.3334
$ws.Range("I32").Value = 1570.5
$ws.Range("K32").Value = 1570.5
$ws.Range("M32").Value = -1283.5
$ws.Range("H61").Value = 2171.9666
$ws.Range("I61").Value = 1204.8462
$ws.Range("J61").Value = 2911.5293
$ws.Range("K61").Value = 1204.8462
$ws.Range("L61").Value = 2911.5293
$ws.Range("M61").Value = -992.8462
$ws.Range("N61").Value = -3335.5293
$ws.Range("H110").Value = 807.8333
$ws.Range("I110").Value = 757
$ws.Range("K110").Value = 757
$ws.Range("M110").Value = 1288
$ws.Range("H116").Value = 2098.6155
$ws.Range("I116").Value = 1688.3
$ws.Range("J116").Value = 3466.3333
$ws.Range("K116").Value = 1688.3
$ws.Range("L116").Value = 3466.3333
$ws.Range("M116").Value = 605.7
$ws.Range("N116").Value = -8054.3333
$ws.Range("H132").Value = 2245.7
$ws.Range("I132").Value = 1996.5667
$ws.Range("J132").Value = 2993.1
$ws.Range("K132").Value = 5989.7001
$ws.Range("L132").Value = 8979.299999999999
$ws.Range("M132").Value = -3459.7001
$ws.Range("N132").Value = -14039.3
$ws.Range("H136").Value = 2171.9666
$ws.Range("I136").Value = 1204.8462
$ws.Range("J136").Value = 2911.5293
$ws.Range("K136").Value = 3614.5386
$ws.Range("L136").Value = 8734.5879
$ws.Range("M136").Value = -1064.5386
$ws.Range("N136").Value = -13834.5879

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2098.6155
$ws.Range("I3").Value = 1688.3
$ws.Range("J3").Value = 3466.3333
$ws.Range("K3").Value = 1688.3
$ws.Range("L3").Value = 3466.3333
$ws.Range("M3").Value = -1574.3
$ws.Range("N3").Value = -3694.3333
$ws.Range("H4").Value = 340.42856
$ws.Range("I4").Value = 314
$ws.Range("K4").Value = 314
$ws.Range("M4").Value = -199
$ws.Range("H22").Value = 796.6667
$ws.Range("I22").Value = 883.875
$ws.Range("J22").Value = 99
$ws.Range("K22").Value = 883.875
$ws.Range("L22").Value = 99
$ws.Range("M22").Value = -710.875
$ws.Range("N22").Value = -445
$ws.Range("H134").Value = 6424.1113
$ws.Range("I134").Value = 6158
$ws.Range("J134").Value = 6756.75
$ws.Range("K134").Value = 18474
$ws.Range("L134").Value = 20270.25
$ws.Range("M134").Value = -15939
$ws.Range("N134").Value = -25340.25

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 457.69232
$ws.Range("I7").Value = 384.1111
$ws.Range("J7").Value = 623.25
$ws.Range("K7").Value = 384.1111
$ws.Range("L7").Value = 623.25
$ws.Range("M7").Value = -271.1111
$ws.Range("N7").Value = -849.25
$ws.Range("H94").Value = 785
$ws.Range("I94").Value = 678
$ws.Range("K94").Value = 678
$ws.Range("M94").Value = -227
$ws.Range("H99").Value = 2880.4443
$ws.Range("I99").Value = 3156.7144
$ws.Range("K99").Value = 3156.7144
$ws.Range("M99").Value = -1658.7144
$ws.Range("H126").Value = 2880.4443
$ws.Range("I126").Value = 3156.7144
$ws.Range("K126").Value = 9470.143199999999
$ws.Range("M126").Value = -7000.143199999999
$ws.Range("H132").Value = 14499954
$ws.Range("I132").Value = 6347.4736
$ws.Range("K132").Value = 19042.4208
$ws.Range("M132").Value = -16512.4208

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I3").Value = 9250
$ws.Range("J3").Value = 7000
$ws.Range("K3").Value = 27750
$ws.Range("L3").Value = 21000
$ws.Range("M3").Value = -27638
$ws.Range("N3").Value = -21224
$ws.Range("H14").Value = 461.625
$ws.Range("I14").Value = 461.625
$ws.Range("K14").Value = 1384.875
$ws.Range("M14").Value = -1211.875
$ws.Range("H113").Value = 1475.1666
$ws.Range("I113").Value = 430
$ws.Range("K113").Value = 1290
$ws.Range("M113").Value = 880

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 618.2
$ws.Range("I55").Value = 279.5
$ws.Range("J55").Value = 844
$ws.Range("K55").Value = 279.5
$ws.Range("L55").Value = 844
$ws.Range("M55").Value = -106.5
$ws.Range("N55").Value = -1190
$ws.Range("H100").Value = 4853.727
$ws.Range("I100").Value = 4321.5557
$ws.Range("J100").Value = 7248.5
$ws.Range("K100").Value = 4321.5557
$ws.Range("L100").Value = 7248.5
$ws.Range("M100").Value = -3780.5557
$ws.Range("N100").Value = -8330.5

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 690.35486
$ws.Range("I107").Value = 594.15
$ws.Range("K107").Value = 1782.45
$ws.Range("M107").Value = 137.5500000000002
$ws.Range("H132").Value = 2270.4187
$ws.Range("I132").Value = 2387.0571
$ws.Range("K132").Value = 7161.1713
$ws.Range("M132").Value = -4631.1713
